# Auto-generated script applying the 2023-06-03 daily crime data update
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 3005
$ws.Range('J3').Value = 3117
$ws.Range('C4').Value = 1824
$ws.Range('I4').Value = 1758
$ws.Range('J4').Value = 702
$ws.Range('J5').Value = 244
$ws.Range('J6').Value = 3742
$ws.Range('C7').Value = 28367
$ws.Range('I7').Value = 26204
$ws.Range('J7').Value = 10810

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 39
$ws.Range('J3').Value = 33
$ws.Range('J7').Value = 121

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J2').Value = 14
$ws.Range('J7').Value = 41

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 108
$ws.Range('J6').Value = 100
$ws.Range('J7').Value = 354

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 53
$ws.Range('J6').Value = 44
$ws.Range('J7').Value = 160

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 83
$ws.Range('J3').Value = 159
$ws.Range('J5').Value = 9
$ws.Range('J7').Value = 389

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 83
$ws.Range('J3').Value = 84
$ws.Range('J6').Value = 101
$ws.Range('J7').Value = 286

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 84
$ws.Range('J5').Value = 31
$ws.Range('J6').Value = 103
$ws.Range('J7').Value = 324
$ws.Range('J8').Value = 695
$ws.Range('J11').Value = 155
$ws.Range('J14').Value = 41
$ws.Range('J15').Value = 125
$ws.Range('J17').Value = 20
$ws.Range('J18').Value = 112
$ws.Range('J19').Value = 345
$ws.Range('J20').Value = 220
$ws.Range('J23').Value = 111
$ws.Range('J26').Value = 17
$ws.Range('J27').Value = 65
$ws.Range('J28').Value = 4
$ws.Range('J29').Value = 614
$ws.Range('J33').Value = 454
$ws.Range('J36').Value = 152
$ws.Range('J37').Value = 354
$ws.Range('J41').Value = 76
$ws.Range('J42').Value = 437
$ws.Range('J43').Value = 95
$ws.Range('J48').Value = 106
$ws.Range('J50').Value = 63
$ws.Range('J52').Value = 292
$ws.Range('J53').Value = 104
$ws.Range('J55').Value = 143
$ws.Range('J60').Value = 71
$ws.Range('C63').Value = 255
$ws.Range('J63').Value = 47
$ws.Range('J65').Value = 286
$ws.Range('J67').Value = 389
$ws.Range('J69').Value = 27
$ws.Range('J76').Value = 155
$ws.Range('J78').Value = 139
$ws.Range('J85').Value = 498
$ws.Range('J86').Value = 63
$ws.Range('I88').Value = 242
$ws.Range('J88').Value = 112
$ws.Range('J89').Value = 121
$ws.Range('J90').Value = 123
$ws.Range('J91').Value = 124
$ws.Range('J94').Value = 95
$ws.Range('J99').Value = 160
$ws.Range('C101').Value = 28367
$ws.Range('I101').Value = 26204
$ws.Range('J101').Value = 10810

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 121
$ws.Range('J3').Value = 141
$ws.Range('J6').Value = 152
$ws.Range('J7').Value = 454

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 179
$ws.Range('J3').Value = 210
$ws.Range('J5').Value = 28
$ws.Range('J6').Value = 160
$ws.Range('J7').Value = 614

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 81
$ws.Range('J3').Value = 97
$ws.Range('J4').Value = 18
$ws.Range('J7').Value = 345

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J4').Value = 18
$ws.Range('J7').Value = 106

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J2').Value = 25
$ws.Range('J4').Value = 16
$ws.Range('J7').Value = 155

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 190
$ws.Range('J7').Value = 498

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J3').Value = 31
$ws.Range('J7').Value = 103

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J2').Value = 19
$ws.Range('J6').Value = 38
$ws.Range('J7').Value = 76

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 93
$ws.Range('J5').Value = 11
$ws.Range('J6').Value = 225
$ws.Range('J7').Value = 437

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J2').Value = 34
$ws.Range('J7').Value = 139

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J2').Value = 36
$ws.Range('J3').Value = 28
$ws.Range('J6').Value = 73
$ws.Range('J7').Value = 143

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J3').Value = 37
$ws.Range('J7').Value = 111

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('J2').Value = 6
$ws.Range('J3').Value = 9
$ws.Range('J7').Value = 27

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J3').Value = 57
$ws.Range('J7').Value = 124

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 73
$ws.Range('J7').Value = 220

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J6').Value = 59
$ws.Range('J7').Value = 112

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('J2').Value = 7
$ws.Range('J7').Value = 20

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J6').Value = 50
$ws.Range('J7').Value = 152

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 67
$ws.Range('J3').Value = 78
$ws.Range('J6').Value = 130
$ws.Range('J7').Value = 292

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J2').Value = 23
$ws.Range('J3').Value = 18
$ws.Range('J7').Value = 95

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J2').Value = 36
$ws.Range('J3').Value = 34
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 125

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J4').Value = 12
$ws.Range('J7').Value = 63

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('J6').Value = 11
$ws.Range('J7').Value = 17

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 56
$ws.Range('J3').Value = 30
$ws.Range('J4').Value = 12
$ws.Range('J7').Value = 155

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J6').Value = 24
$ws.Range('J7').Value = 84

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J2').Value = 26
$ws.Range('J3').Value = 35
$ws.Range('I4').Value = 5
$ws.Range('J4').Value = 2
$ws.Range('I7').Value = 242
$ws.Range('J7').Value = 112

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 211
$ws.Range('J3').Value = 222
$ws.Range('J5').Value = 21
$ws.Range('J6').Value = 205
$ws.Range('J7').Value = 695

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('J3').Value = 3
$ws.Range('J7').Value = 31

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J3').Value = 16
$ws.Range('J4').Value = 5
$ws.Range('J7').Value = 65

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J4').Value = 30
$ws.Range('J7').Value = 63

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J3').Value = 38
$ws.Range('J7').Value = 123

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('J2').Value = 23
$ws.Range('J6').Value = 23
$ws.Range('J7').Value = 71

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('J6').Value = 60
$ws.Range('J7').Value = 95

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J6').Value = 62
$ws.Range('J7').Value = 104

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 94
$ws.Range('J7').Value = 324

$ws = $wb.Worksheets.Item('Edison Park')
$ws.Range('J4').Value = 1
$ws.Range('J7').Value = 4
